$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 78

# New Pick 3 draw result row, appended by the automated results updater.
# Date and Phase look numeric to Excel's type-inference, so they are written
# with a leading apostrophe to force literal text (matching the rest of the
# "numberStoredAsText" column data above them). Game/Result/InsertedAt are
# plain text already and don't need that treatment.
$ws.Cells.Item($row, 1).Formula = "'2025-12-03"
$ws.Cells.Item($row, 2).Value   = "Pick 3"
$ws.Cells.Item($row, 3).Formula = "'251203"
$ws.Cells.Item($row, 4).Value   = "2-8-5"
$ws.Cells.Item($row, 5).Value   = "2025-12-03T21:43:55.285+04:00"
